# Daily attendance processing - 2026-01-25 04:09:33
# Swap the order of "System" and the email address in column G
# ("System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System")
# for every row in the "Recorded By" column where that exact text occurs.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldText = "System, dnasr281@gmail.com"
$newText = "dnasr281@gmail.com, System"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($row = 1; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)   # Column G
    $val = $cell.Value2
    if ($val -eq $oldText) {
        $cell.Value = $newText
    }
}
